$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) whose values change for rows 2-25.
$cols = @("B","C","D","E","I","J","K","L","N")

# Each entry: row number followed by the new values for columns B,C,D,E,I,J,K,L,N (in that order)
$data = @(
    @(2,1.02,1.030226741290073,1.033412852184459,1.029905282718969,1.031716128618642,1.035369154157195,1.036215274259204,1.0327178385106,1.036839497114612),
    @(3,1.02,1.031775028515915,1.034589960293003,1.03123895177173,1.032092953585462,1.036555751098005,1.037200374786341,1.033858338331424,1.038027779159171),
    @(4,1.02,1.032774568074945,1.035349297015272,1.03210015501105,1.032334012074049,1.03732093072665,1.037834908388872,1.03459402218868,1.038794045430656),
    @(5,1.02,1.03319423222559,1.035667970444238,1.03246178866055,1.032434691966151,1.037641990321572,1.038100979742855,1.034902760959787,1.039115560967029),
    @(6,1.02,1.033264664164893,1.035721444943397,1.03252248435277,1.032451557871977,1.037695861378855,1.038145614201886,1.034954567931677,1.039169508527392),
    @(7,1.02,1.032780177765127,1.035353557304281,1.032104988802436,1.03233535995761,1.037325223174315,1.037838466336867,1.034598149694148,1.038798343974089),
    @(8,1.02,1.030750475522451,1.033811147828481,1.030356373407867,1.031844054034412,1.035770718641579,1.036548796007329,1.033103755290148,1.037241631866612),
    @(9,1.02,1.027155782371458,1.031075088202338,1.027261196589756,1.030956956805027,1.033011030626097,1.03425381581958,1.030452547496526,1.034478024777608),
    @(10,1.02,1.024746534676024,1.029238451427729,1.025187927933811,1.030351031747408,1.031157028538979,1.032708383452025,1.028672606827912,1.032621389794869),
    @(11,1.02,1.023700132845546,1.028440088702525,1.024287744137892,1.030085175253051,1.030350751470098,1.032035446240189,1.027898817046721,1.031813967720068),
    @(12,1.02,1.023310962568387,1.028143069508344,1.023952999172307,1.029985897008691,1.030050732130064,1.031784915874315,1.027610928198857,1.031513522318159),
    @(13,1.02,1.023394463284066,1.028206802609395,1.024024820297412,1.030007216446351,1.030115111523073,1.031838681470702,1.027672702690045,1.031577993137291),
    @(14,1.02,1.023667973975892,1.028415546672843,1.024260081751877,1.03007697965928,1.030325962685809,1.032014749024051,1.027875029679511,1.031789143732863),
    @(15,1.02,1.023836427642071,1.028544098031022,1.024404983911852,1.030119893100763,1.030455804261363,1.032123154099202,1.027999627627299,1.031919169698347),
    @(16,1.02,1.024815912860513,1.029291370417966,1.025247617865366,1.030368602009631,1.031210464365005,1.03275296433902,1.028723895427781,1.032674901505897),
    @(17,1.02,1.025429457776522,1.029759282606088,1.02577551879258,1.03052367457522,1.031682904163796,1.033147017195313,1.029177383208201,1.033148012223392),
    @(18,1.02,1.025787021979408,1.030031910573024,1.026083199220714,1.030613789573994,1.031958134964244,1.033376499781454,1.029441599937607,1.033423633883146),
    @(19,1.02,1.025908890563228,1.030124819477024,1.026188070673161,1.030644459536498,1.032051924945977,1.033454686271623,1.029531641231789,1.033517557057411),
    @(20,1.02,1.025363661968792,1.029709110864879,1.025718904419673,1.030507071556604,1.031632250595618,1.033104776527076,1.029128758838539,1.033097286721337),
    @(21,1.02,1.023587445519162,1.028354089835456,1.024190813589252,1.030056450715858,1.0302638870362,1.031962917354706,1.027815462436439,1.031726979928711),
    @(22,1.02,1.022467826431985,1.027499401031799,1.023227858522405,1.029770074950453,1.029400459236897,1.031241674784861,1.026987023811182,1.030862325962897),
    @(23,1.02,1.023061631066566,1.02795274950535,1.023738549350206,1.029922178666514,1.029858473948305,1.031624335399298,1.027426455233388,1.03132099110773),
    @(24,1.02,1.025393393208265,1.02973178223391,1.025744486744897,1.030514574784977,1.031655139799013,1.033123864384824,1.029150731009874,1.033120208430026),
    @(25,1.02,1.028087304642597,1.031784615693305,1.02806306915442,1.031188840549312,1.033726946793717,1.034849817890729,1.031140113594051,1.035194957628302)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
